# Add the "I0" and "IF" columns (I and J) to the pitching log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (copy H1's header formatting - bold, bordered, centered -
# onto the two new header cells)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2..84: pairs of (I0, IF) values taken directly from the
# source data (not derivable from the other columns).
$data = @(
    @(4,4),   # row 2
    @(9,9),   # row 3
    @(8,8),   # row 4
    @(9,9),   # row 5
    @(1,1),   # row 6
    @(5,5),   # row 7
    @(8,8),   # row 8
    @(3,3),   # row 9
    @(8,8),   # row 10
    @(8,8),   # row 11
    @(8,8),   # row 12
    @(9,9),   # row 13
    @(9,9),   # row 14
    @(9,9),   # row 15
    @(8,8),   # row 16
    @(8,8),   # row 17
    @(9,9),   # row 18
    @(8,8),   # row 19
    @(7,7),   # row 20
    @(9,9),   # row 21
    @(6,6),   # row 22
    @(9,9),   # row 23
    @(9,9),   # row 24
    @(8,8),   # row 25
    @(9,9),   # row 26
    @(9,9),   # row 27
    @(9,9),   # row 28
    @(6,6),   # row 29
    @(8,8),   # row 30
    @(7,8),   # row 31
    @(7,7),   # row 32
    @(7,7),   # row 33
    @(7,7),   # row 34
    @(5,5),   # row 35
    @(8,8),   # row 36
    @(7,8),   # row 37
    @(9,9),   # row 38
    @(9,9),   # row 39
    @(10,10), # row 40
    @(9,9),   # row 41
    @(8,9),   # row 42
    @(9,9),   # row 43
    @(9,9),   # row 44
    @(9,9),   # row 45
    @(9,9),   # row 46
    @(9,9),   # row 47
    @(9,9),   # row 48
    @(8,8),   # row 49
    @(9,9),   # row 50
    @(9,9),   # row 51
    @(8,8),   # row 52
    @(9,9),   # row 53
    @(9,9),   # row 54
    @(9,10),  # row 55
    @(9,9),   # row 56
    @(8,8),   # row 57
    @(8,8),   # row 58
    @(9,9),   # row 59
    @(9,9),   # row 60
    @(9,9),   # row 61
    @(8,9),   # row 62
    @(9,9),   # row 63
    @(9,9),   # row 64
    @(9,9),   # row 65
    @(8,8),   # row 66
    @(9,9),   # row 67
    @(8,8),   # row 68
    @(8,8),   # row 69
    @(8,8),   # row 70
    @(8,8),   # row 71
    @(9,9),   # row 72
    @(9,9),   # row 73
    @(8,9),   # row 74
    @(8,8),   # row 75
    @(8,9),   # row 76
    @(7,8),   # row 77
    @(8,8),   # row 78
    @(5,5),   # row 79
    @(6,6),   # row 80
    @(4,4),   # row 81
    @(1,3),   # row 82
    @(1,3),   # row 83
    @(1,2)    # row 84
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
    $row++
}
